$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43-51 down to 44-52.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new data record.
$ws.Cells.Item(43, 1).Value = 2
$ws.Cells.Item(43, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 44588
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(43, 6).Value = 100112032
$ws.Cells.Item(43, 7).Value = "Zapallo italiano"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 300
$ws.Cells.Item(43, 11).Value = 10000
$ws.Cells.Item(43, 12).Value = 11000
$ws.Cells.Item(43, 13).Value = 10500
$ws.Cells.Item(43, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 175
$ws.Cells.Item(43, 17).Value = 60
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Ensure the date style/number format (used by column D on adjacent rows) is preserved.
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
